$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "Price" (column D) cells whose new values would otherwise be
# auto-detected as numbers by Excel. Each one is written as a literal-text
# formula, then converted in place to a plain value via PasteSpecial so the
# stored cell keeps its original (unstyled, text) representation.
$c = $ws.Range("D5")
$c.Formula = '="587.23"'
$c.Copy()
$c.PasteSpecial(-4163)

$c = $ws.Range("D6")
$c.Formula = '="134.28"'
$c.Copy()
$c.PasteSpecial(-4163)

$c = $ws.Range("D10")
$c.Formula = '="7.26"'
$c.Copy()
$c.PasteSpecial(-4163)

$c = $ws.Range("D17")
$c.Formula = '="25.67"'
$c.Copy()
$c.PasteSpecial(-4163)

$c = $ws.Range("D19")
$c.Formula = '="5.75"'
$c.Copy()
$c.PasteSpecial(-4163)

$c = $ws.Range("D20")
$c.Formula = '="13.55"'
$c.Copy()
$c.PasteSpecial(-4163)

$c = $ws.Range("D21")
$c.Formula = '="394.19"'
$c.Copy()
$c.PasteSpecial(-4163)

$c = $ws.Range("D24")
$c.Formula = '="74.69"'
$c.Copy()
$c.PasteSpecial(-4163)

$c = $ws.Range("D28")
$c.Formula = '="0.998"'
$c.Copy()
$c.PasteSpecial(-4163)

$c = $ws.Range("D31")
$c.Formula = '="8.22"'
$c.Copy()
$c.PasteSpecial(-4163)

$c = $ws.Range("D32")
$c.Formula = '="1.48"'
$c.Copy()
$c.PasteSpecial(-4163)

$c = $ws.Range("D36")
$c.Formula = '="23.41"'
$c.Copy()
$c.PasteSpecial(-4163)

$c = $ws.Range("D37")
$c.Formula = '="5.14"'
$c.Copy()
$c.PasteSpecial(-4163)

$c = $ws.Range("D40")
$c.Formula = '="166.80"'
$c.Copy()
$c.PasteSpecial(-4163)

$c = $ws.Range("D41")
$c.Formula = '="0.0778"'
$c.Copy()
$c.PasteSpecial(-4163)

$c = $ws.Range("D44")
$c.Formula = '="25.42"'
$c.Copy()
$c.PasteSpecial(-4163)

# Remaining "Price" and "Volume(1h)" text cells can be written directly —
# their new values are not valid numeric literals, so Excel keeps them as text.
$ws.Range("D2").Value = "64.257.94"
$ws.Range("E2").Value = "  +0.20%  "
$ws.Range("D3").Value = "3.494.69"
$ws.Range("E3").Value = "  -0.71%  "
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("E5").Value = "  +0.23%  "
$ws.Range("E6").Value = "  +0.87%  "
$ws.Range("E7").Value = "  +0.23%  "
$ws.Range("E8").Value = "  +0.07%  "
$ws.Range("E9").Value = "  -0.31%  "
$ws.Range("E10").Value = "  +1.95%  "
$ws.Range("E11").Value = "  +1.95%  "
$ws.Range("D12").Value = "4.088.22"
$ws.Range("E12").Value = "  -0.84%  "
$ws.Range("E13").Value = "  +1.09%  "
$ws.Range("E14").Value = "  +1.36%  "
$ws.Range("D15").Value = "3.492.31"
$ws.Range("E15").Value = "  -0.84%  "
$ws.Range("D16").Value = "64.338.32"
$ws.Range("E17").Value = "  -6.97%  "
$ws.Range("E18").Value = "  +1.00%  "
$ws.Range("E19").Value = "  +2.18%  "
$ws.Range("E20").Value = "  -3.26%  "
$ws.Range("E21").Value = "  +2.31%  "
$ws.Range("E22").Value = "  -0.83%  "
$ws.Range("D23").Value = "3.633.79"
$ws.Range("E23").Value = "  -0.78%  "
$ws.Range("E24").Value = "  +0.96%  "
$ws.Range("E25").Value = "  +0.07%  "
$ws.Range("E26").Value = "  +0.15%  "
$ws.Range("E27").Value = "  +0.50%  "
$ws.Range("E28").Value = "  +0.34%  "
$ws.Range("E29").Value = "  -1.87%  "
$ws.Range("E30").Value = "  +0.02%  "
$ws.Range("E31").Value = "  -2.41%  "
$ws.Range("E32").Value = "  -5.80%  "
$ws.Range("D33").Value = "3.516.02"
$ws.Range("E33").Value = "  -0.40%  "
$ws.Range("E34").Value = "  +3.97%  "
$ws.Range("E36").Value = "  -0.63%  "
$ws.Range("E37").Value = "  -4.14%  "
$ws.Range("E38").Value = "  -0.51%  "
$ws.Range("E39").Value = "  -0.59%  "
$ws.Range("E40").Value = "  +3.81%  "
$ws.Range("E41").Value = "  -1.31%  "
$ws.Range("E42").Value = "  -1.04%  "
$ws.Range("E43").Value = "  -0.11%  "
$ws.Range("E44").Value = "  -4.75%  "
$ws.Range("E45").Value = "  -0.77%  "
$ws.Range("E46").Value = "  +2.10%  "
$ws.Range("E47").Value = "  -3.87%  "
$ws.Range("D48").Value = "2.462.71"
$ws.Range("E48").Value = "  -0.42%  "
$ws.Range("E49").Value = "  -0.95%  "
$ws.Range("E50").Value = "  -1.54%  "
$ws.Range("E51").Value = "  -1.13%  "
